# Update sample-plant assignments and revenue totals per the source data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "SVEKLA"
$ws.Range("D2").Value = "SOYA"
$ws.Range("E2").Value = "PAR"
$ws.Range("C3").Value = "PSHENICA"
$ws.Range("E3").Value = "PAR"
$ws.Range("C4").Value = "PSHENICA"
$ws.Range("D4").Value = "PSHENICA"
$ws.Range("C5").Value = "PSHENICA"
$ws.Range("D5").Value = "YACHMEN"
$ws.Range("E5").Value = "PAR"
$ws.Range("C6").Value = "YACHMEN"
$ws.Range("D6").Value = "YACHMEN"
$ws.Range("E6").Value = "SVEKLA"
$ws.Range("D7").Value = "SVEKLA"
$ws.Range("C8").Value = "SOYA"
$ws.Range("D8").Value = "PAR"
$ws.Range("E8").Value = "PSHENICA"
$ws.Range("C9").Value = "SOYA"
$ws.Range("D9").Value = "YACHMEN"
$ws.Range("E9").Value = "PSHENICA"
$ws.Range("D10").Value = "SOYA"
$ws.Range("E11").Value = "PAR"
$ws.Range("C12").Value = "YACHMEN"
$ws.Range("D12").Value = "YACHMEN"
$ws.Range("E12").Value = "SOYA"
$ws.Range("C13").Value = "PSHENICA"
$ws.Range("D13").Value = "SOYA"
$ws.Range("E13").Value = "PAR"
$ws.Range("D14").Value = "SOYA"
$ws.Range("E14").Value = "YACHMEN"
$ws.Range("C15").Value = "PAR"
$ws.Range("D15").Value = "YACHMEN"
$ws.Range("E15").Value = "SVEKLA"
$ws.Range("C16").Value = "SVEKLA"
$ws.Range("E16").Value = "PAR"
$ws.Range("C17").Value = "YACHMEN"
$ws.Range("D17").Value = "PAR"
$ws.Range("E17").Value = "YACHMEN"
$ws.Range("C18").Value = "PAR"
$ws.Range("D18").Value = "PAR"
$ws.Range("E18").Value = "YACHMEN"
$ws.Range("D19").Value = "PAR"
$ws.Range("E19").Value = "SVEKLA"
$ws.Range("C20").Value = "SVEKLA"
$ws.Range("D20").Value = "SVEKLA"
$ws.Range("C21").Value = "SOYA"
$ws.Range("E21").Value = "PSHENICA"
$ws.Range("C22").Value = "PSHENICA"
$ws.Range("D22").Value = "PSHENICA"
$ws.Range("E22").Value = "PSHENICA"
$ws.Range("C23").Value = "PAR"
$ws.Range("D23").Value = "SVEKLA"
$ws.Range("E23").Value = "PSHENICA"
$ws.Range("C24").Value = "SOYA"
$ws.Range("D24").Value = "SVEKLA"
$ws.Range("E24").Value = "YACHMEN"
$ws.Range("C25").Value = "PSHENICA"
$ws.Range("D25").Value = "YACHMEN"
$ws.Range("E25").Value = "YACHMEN"
$ws.Range("D26").Value = "PAR"
$ws.Range("E26").Value = "PAR"
$ws.Range("C27").Value = "PAR"
$ws.Range("D27").Value = "PSHENICA"
$ws.Range("E27").Value = "YACHMEN"
$ws.Range("C28").Value = "SOYA"
$ws.Range("E28").Value = "SVEKLA"
$ws.Range("C29").Value = "PAR"
$ws.Range("D29").Value = "YACHMEN"
$ws.Range("E29").Value = "YACHMEN"
$ws.Range("E30").Value = "PAR"
$ws.Range("C31").Value = "PAR"
$ws.Range("E31").Value = "SOYA"
$ws.Range("C32").Value = "SVEKLA"
$ws.Range("E33").Value = "SOYA"
$ws.Range("C34").Value = "YACHMEN"
$ws.Range("D34").Value = "PAR"
$ws.Range("C35").Value = "PAR"
$ws.Range("D35").Value = "YACHMEN"
$ws.Range("E35").Value = "PAR"
$ws.Range("C36").Value = "SOYA"
$ws.Range("D36").Value = "SVEKLA"
$ws.Range("E36").Value = "PSHENICA"
$ws.Range("E37").Value = "SOYA"
$ws.Range("D38").Value = "PAR"
$ws.Range("C39").Value = "SVEKLA"
$ws.Range("D39").Value = "PSHENICA"
$ws.Range("E39").Value = "SOYA"
$ws.Range("C40").Value = "PAR"
$ws.Range("D40").Value = "PSHENICA"
$ws.Range("E40").Value = "PSHENICA"
$ws.Range("C41").Value = "YACHMEN"
$ws.Range("E41").Value = "YACHMEN"
$ws.Range("D42").Value = "YACHMEN"
$ws.Range("E42").Value = "SOYA"
$ws.Range("E43").Value = "YACHMEN"
$ws.Range("D44").Value = "SVEKLA"
$ws.Range("E44").Value = "PSHENICA"
$ws.Range("C45").Value = "PSHENICA"
$ws.Range("D45").Value = "SOYA"
$ws.Range("E45").Value = "SOYA"
$ws.Range("D46").Value = "PAR"
$ws.Range("E46").Value = "PSHENICA"
$ws.Range("C47").Value = "SVEKLA"
$ws.Range("D47").Value = "YACHMEN"
$ws.Range("E47").Value = "SOYA"
$ws.Range("C48").Value = "PAR"
$ws.Range("D48").Value = "PSHENICA"
$ws.Range("C49").Value = "YACHMEN"
$ws.Range("D49").Value = "PSHENICA"
$ws.Range("E49").Value = "PSHENICA"
$ws.Range("C50").Value = "PAR"
$ws.Range("E50").Value = "SVEKLA"
$ws.Range("D51").Value = "PSHENICA"
$ws.Range("E51").Value = "PAR"
$ws.Range("C52").Value = "YACHMEN"
$ws.Range("D52").Value = "SVEKLA"
$ws.Range("E52").Value = "YACHMEN"
$ws.Range("C53").Value = "PSHENICA"
$ws.Range("D53").Value = "YACHMEN"
$ws.Range("E53").Value = "PAR"
$ws.Range("C54").Value = "PSHENICA"
$ws.Range("D54").Value = "SVEKLA"
$ws.Range("E54").Value = "PSHENICA"
$ws.Range("C55").Value = "SOYA"
$ws.Range("C56").Value = "SOYA"
$ws.Range("D56").Value = "SOYA"
$ws.Range("E56").Value = "PSHENICA"
$ws.Range("C57").Value = "YACHMEN"
$ws.Range("E57").Value = "PSHENICA"
$ws.Range("C58").Value = "SVEKLA"
$ws.Range("D58").Value = "PAR"
$ws.Range("E58").Value = "PSHENICA"
$ws.Range("C59").Value = "PSHENICA"
$ws.Range("D59").Value = "SOYA"
$ws.Range("E59").Value = "SVEKLA"
$ws.Range("C60").Value = "SOYA"
$ws.Range("D60").Value = "YACHMEN"
$ws.Range("E60").Value = "YACHMEN"
$ws.Range("C61").Value = 436149608
$ws.Range("D61").Value = 420653496
$ws.Range("E61").Value = 417776608
